$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update n_requests (E), temp (F), save (G) columns for data rows 2-37
$ws.Range("E2:E37").Value = 3
$ws.Range("F2:F37").Value = 0.5
$ws.Range("G2:G37").Value = 1

# Update the active selection shown in the sheet view
$ws.Activate()
$ws.Range("H5").Select()
